# Weekly price-sheet update: a new weekly price entry is inserted as row 48
# (pushing the existing rows 48:160 down to 49:161); every column value for
# the new row mirrors the sheet's constant columns (market/region/product
# metadata) while the variable columns (date, volume, prices) carry the
# new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 48; rows 48:160 shift down to 49:161.
$ws.Rows(48).Insert()

$row = 48

$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = (Get-Date -Year 2021 -Month 11 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = 100112005
$ws.Cells.Item($row, 7).Value = "Puerro"
$ws.Cells.Item($row, 8).Value = "Azul de Maquehue"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 70
$ws.Cells.Item($row, 11).Value = 8000
$ws.Cells.Item($row, 12).Value = 8000
$ws.Cells.Item($row, 13).Value = 8000
$ws.Cells.Item($row, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item($row, 15).Value = "Provincia de Cautín"
$ws.Cells.Item($row, 16).Value = 667
$ws.Cells.Item($row, 17).Value = 12
$ws.Cells.Item($row, 18).Value = "Hortaliza"
